$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: "CS4102, Spring 2021" -> "CS4102, Fall 2021"
$para1 = $tr.Paragraphs(1)
$para1.Runs(1).Text = "CS4102, Fall 2021"

# Paragraph 3: merge the "Readings" run into the ": CLRS 23.2, 24.2, 24.3" run
$para3 = $tr.Paragraphs(3)
$para3.Runs(1).Text = ""
$para3 = $tr.Paragraphs(3)
$para3.Runs(1).Text = "Readings: CLRS 23.2, 24.2, 24.3"
